# Update "想去人数" (want-to-go count) figures for the 南宁-漫展信息 workbook
# to match the data refresh output at commit 456a3b4.
#
# Sheet "展览" (Exhibitions): rows 2-4, column F
#   F2: 970  -> 975
#   F3: 1941 -> 1964
#   F4: 428  -> 433
#
# Sheet "全部类型" (All types): rows 4-6, column F (same events, same new counts)
#   F4: 970  -> 975
#   F5: 1941 -> 1964
#   F6: 428  -> 433

$wb = $excel.ActiveWorkbook

$wsExhibition = $wb.Worksheets.Item("展览")
$wsExhibition.Range("F2").Value = 975
$wsExhibition.Range("F3").Value = 1964
$wsExhibition.Range("F4").Value = 433

$wsAllTypes = $wb.Worksheets.Item("全部类型")
$wsAllTypes.Range("F4").Value = 975
$wsAllTypes.Range("F5").Value = 1964
$wsAllTypes.Range("F6").Value = 433
